$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.309.53"
$ws.Range("E2").Value = "  -0.80%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.706.42"
$ws.Range("E3").Value = "  -1.06%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.86"
$ws.Range("E5").Value = "  -0.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5322"
$ws.Range("E6").Value = "  -0.94%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2666"
$ws.Range("E8").Value = "  -0.38%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06601"
$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.81"
$ws.Range("E10").Value = "  -4.32%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07649"
$ws.Range("E11").Value = "  -1.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.514"
$ws.Range("E12").Value = "  -2.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.940.42"
$ws.Range("E13").Value = "  -1.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.690.22"
$ws.Range("E14").Value = "  -1.96%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5824"
$ws.Range("E15").Value = "  -0.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8176"
$ws.Range("E16").Value = "  -1.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.62"
$ws.Range("E17").Value = "  -0.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "27.324.33"
$ws.Range("E18").Value = "  -0.82%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.44"
$ws.Range("E19").Value = "  -2.12%  "

$ws.Range("E20").Value = "  -0.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.636"
$ws.Range("E21").Value = "  -2.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.41"
$ws.Range("E22").Value = "  -2.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.984"
$ws.Range("E23").Value = "  -1.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.004"
$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.57"
$ws.Range("E25").Value = "  -3.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.707"
$ws.Range("E26").Value = "  +0.67%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1203"
$ws.Range("E27").Value = "  -2.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.223"
$ws.Range("E28").Value = "  -2.38%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.18"
$ws.Range("E29").Value = "  -2.85%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05375"
$ws.Range("E30").Value = "  -2.98%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.288"
$ws.Range("E31").Value = "  -1.17%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.475"
$ws.Range("E32").Value = "  -2.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.410"
$ws.Range("E33").Value = "  -1.48%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.646"
$ws.Range("E34").Value = "  -1.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.866"
$ws.Range("E35").Value = "  +1.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9496"
$ws.Range("E36").Value = "  -1.21%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.406"
$ws.Range("E37").Value = "  -1.61%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5857"
$ws.Range("E38").Value = "  -1.61%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01639"
$ws.Range("E39").Value = "  -0.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.807"
$ws.Range("E40").Value = "  -2.14%  "

$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.043.68"
$ws.Range("E41").Value = "  -1.36%  "

$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.003"
$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8425"
$ws.Range("E43").Value = "  -1.43%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.89"
$ws.Range("E44").Value = "  -0.63%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.849.07"
$ws.Range("E45").Value = "  -1.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈110"
$ws.Range("E46").Value = "  -3.91%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.90"
$ws.Range("E47").Value = "  -1.95%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4524"
$ws.Range("E48").Value = "  +1.86%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.005"
$ws.Range("E49").Value = "  +0.39%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.073"
$ws.Range("E50").Value = "  -1.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05232"
$ws.Range("E51").Value = "  -0.83%  "
